$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New case rows (473-488) for case 21TRD09386 / Hemmeter.
# Columns D (statute #) and H/I (dollar amounts) look numeric as plain
# strings, so pre-format those ranges as Text to keep them literal text
# (matching the source data) instead of being coerced into numbers.
$ws.Range("D473:D488").NumberFormat = "@"
$ws.Range("H473:I488").NumberFormat = "@"

$rows = @(
    @("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","DUS UCM","4510.111","UCM","No Contest","Guilty","$ 0","$ 0","None","None"),
    @("21TRD09386","Hemmeter","TAIL LIGHTS-REAR LICENSE PLATE","4513.05","MM","No Contest","Guilty","$ 0","$ 0","None","None")
)

$startRow = 473
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

Write-Host "Rows written: $($rows.Count)"
